$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.83"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'24.25"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'5.370"
$ws.Range("G4").Value = "'12"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'6.521"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'3.148"
$ws.Range("G7").Value = "'12"
$ws.Range("D8").Value = "'0.8168"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.8678"
$ws.Range("G9").Value = "'12"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'0.06991"
$ws.Range("G11").Value = "'12"
$ws.Range("D12").Value = "'0.03270"
$ws.Range("G12").Value = "'12"
$ws.Range("D13").Value = "'0.02898"
$ws.Range("G13").Value = "'12"
$ws.Range("D14").Value = "'0.09369"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'3.732"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.001532"
$ws.Range("G16").Value = "'12"
$ws.Range("D17").Value = "'0.04711"
$ws.Range("G17").Value = "'12"
$ws.Range("D18").Value = "'0.0005978"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'0.006199"
$ws.Range("G19").Value = "'12"
$ws.Range("D20").Value = "'0.001243"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.003844"
$ws.Range("G21").Value = "'12"
$ws.Range("D22").Value = "'0.00008799"
$ws.Range("G22").Value = "'12"
$ws.Range("D23").Value = "'3.531"
$ws.Range("G23").Value = "'12"
$ws.Range("D24").Value = "'2.149"
$ws.Range("G24").Value = "'12"
$ws.Range("D25").Value = "'0.3154"
$ws.Range("G25").Value = "'12"
$ws.Range("D26").Value = "'0.1330"
$ws.Range("G26").Value = "'12"
$ws.Range("D27").Value = "'0.1328"
$ws.Range("G27").Value = "'12"
$ws.Range("D28").Value = "'0.0003013"
$ws.Range("E28").Value = "27UpBotsUBXTBestin24h"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.03702"
$ws.Range("G40").Value = "'12"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1056"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("G41").Value = "'12"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002211"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").Value = "'12"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003068"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.008656"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005381"
$ws.Range("G45").Value = "'12"
$ws.Range("G46").Value = "'12"
$ws.Range("D47").Value = "'0.3883"
$ws.Range("G47").Value = "'12"
$ws.Range("D48").Value = "'0.002565"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("G48").Value = "'12"
$ws.Range("G49").Value = "'12"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"
